# StudySchedule3: add a "day of week" column (D) next to the existing
# study-pair schedule in column B, and highlight the first 12 rows
# (first ~1.5 weeks) in yellow to mark them as done/reviewed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 of the schedule corresponds to a Saturday, then the week repeats.
$days = @("sat", "sun", "mon", "tue", "wed", "thu", "fri")

for ($i = 1; $i -le 28; $i++) {
    $dayName = $days[($i - 1) % 7]
    $ws.Cells.Item($i, 4).Value = $dayName
}

# Highlight the first 12 rows (columns A-D) in yellow.
$ws.Range("A1:D12").Interior.Color = 65535

# Leave the selection on B13, matching where editing left off.
$ws.Range("B13").Select()
